# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 11, pushing the existing
# rows 11-17 down to rows 12-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11 (shifts rows 11-17 -> 12-18)
$ws.Rows.Item(11).Insert()

# Copy static/common values from row 12 (the row that used to be row 11
# before the insert) into the new row 11, since every record in this
# block shares the same market/product metadata.
$ws.Range("A12:T12").Copy()
$ws.Range("A11:T11").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# Now overwrite the values specific to this new weekly record.
$ws.Cells.Item(11, 4).Value = 44917    # D11 Fecha
$ws.Cells.Item(11, 12).Value = "Primera"  # L11 Calidad
$ws.Cells.Item(11, 13).Value = 200     # M11 Volumen
$ws.Cells.Item(11, 14).Value = 7000    # N11 Precio minimo
$ws.Cells.Item(11, 15).Value = 7500    # O11 Precio maximo
$ws.Cells.Item(11, 16).Value = 7250    # P11 Precio promedio ponderado
$ws.Cells.Item(11, 19).Value = 3625    # S11 Precio $/Kg

$wb.Save()
